# Atualização de bases das ligas, do dia: 11-04-2024 às 23:56
#
# This league export had several fixture rows whose stat columns (id,
# score, odds, Asian-handicap lines, P/L, ...) were transposed between
# two adjacent rows, and one not-yet-played fixture (match 8021846,
# Real Santa Cruz vs The Strongest) that should not have been included
# got removed from the bottom of the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 27: restore the correct fixture data for this row
$ws.Range("B27").Value = 6504313
$ws.Range("F27").Value = 'Guabira'
$ws.Range("G27").Value = 'Atletico Palmaflor Vinto'
$ws.Range("H27").Value = 1
$ws.Range("I27").Value = 0
$ws.Range("J27").Value = 'H'
$ws.Range("K27").Value = 1.75
$ws.Range("L27").Value = 3.6
$ws.Range("M27").Value = 4
$ws.Range("N27").Value = 2
$ws.Range("O27").Value = 3.6
$ws.Range("P27").Value = 3.6
$ws.Range("Q27").Value = -0.25
$ws.Range("R27").Value = 1.75
$ws.Range("S27").Value = 2.05
$ws.Range("T27").Value = 2.5
$ws.Range("U27").Value = 1.85
$ws.Range("V27").Value = 1.95
$ws.Range("W27").Value = 1
$ws.Range("X27").Value = -1
$ws.Range("Y27").Value = -1
$ws.Range("Z27").Value = 0.75
$ws.Range("AA27").Value = -1
$ws.Range("AB27").Value = -1
$ws.Range("AC27").Value = 0.95

# Row 28: restore the correct fixture data for this row
$ws.Range("B28").Value = 6504831
$ws.Range("F28").Value = 'Libertad Gran Mamore FC'
$ws.Range("G28").Value = 'Always Ready'
$ws.Range("H28").Value = 1
$ws.Range("I28").Value = 2
$ws.Range("J28").Value = 'A'
$ws.Range("K28").Value = 3.3
$ws.Range("L28").Value = 3.6
$ws.Range("M28").Value = 1.909
$ws.Range("N28").Value = 3.1
$ws.Range("O28").Value = 3.5
$ws.Range("P28").Value = 2.2
$ws.Range("Q28").Value = 0.25
$ws.Range("R28").Value = 1.9
$ws.Range("S28").Value = 1.9
$ws.Range("T28").Value = 2.75
$ws.Range("U28").Value = 2.025
$ws.Range("V28").Value = 1.775
$ws.Range("W28").Value = -1
$ws.Range("X28").Value = -1
$ws.Range("Y28").Value = 1.2
$ws.Range("Z28").Value = -1
$ws.Range("AA28").Value = 0.8999999999999999
$ws.Range("AB28").Value = 0.5125
$ws.Range("AC28").Value = -0.5

# Row 47: restore the correct fixture data for this row
$ws.Range("B47").Value = 6504291
$ws.Range("F47").Value = 'Club Aurora'
$ws.Range("G47").Value = 'Universitario De Vinto'
$ws.Range("H47").Value = 0
$ws.Range("I47").Value = 1
$ws.Range("J47").Value = 'A'
$ws.Range("K47").Value = 2
$ws.Range("L47").Value = 3.2
$ws.Range("M47").Value = 3.3
$ws.Range("N47").Value = 1.7
$ws.Range("O47").Value = 3.6
$ws.Range("P47").Value = 5.5
$ws.Range("Q47").Value = -0.75
$ws.Range("R47").Value = 1.85
$ws.Range("S47").Value = 1.95
$ws.Range("T47").Value = 2.5
$ws.Range("U47").Value = 1.9
$ws.Range("V47").Value = 1.9
$ws.Range("W47").Value = -1
$ws.Range("X47").Value = -1
$ws.Range("Y47").Value = 4.5
$ws.Range("Z47").Value = -1
$ws.Range("AA47").Value = 0.95
$ws.Range("AB47").Value = -1
$ws.Range("AC47").Value = 0.8999999999999999

# Row 48: restore the correct fixture data for this row
$ws.Range("B48").Value = 6504578
$ws.Range("F48").Value = 'Guabira'
$ws.Range("G48").Value = 'Real Santa Cruz'
$ws.Range("H48").Value = 1
$ws.Range("I48").Value = 2
$ws.Range("J48").Value = 'A'
$ws.Range("K48").Value = 2
$ws.Range("L48").Value = 3.2
$ws.Range("M48").Value = 3.4
$ws.Range("N48").Value = 1.75
$ws.Range("O48").Value = 3.6
$ws.Range("P48").Value = 5.25
$ws.Range("Q48").Value = -0.75
$ws.Range("R48").Value = 1.95
$ws.Range("S48").Value = 1.85
$ws.Range("T48").Value = 2.5
$ws.Range("U48").Value = 2
$ws.Range("V48").Value = 1.8
$ws.Range("W48").Value = -1
$ws.Range("X48").Value = -1
$ws.Range("Y48").Value = 4.25
$ws.Range("Z48").Value = -1
$ws.Range("AA48").Value = 0.8500000000000001
$ws.Range("AB48").Value = 1
$ws.Range("AC48").Value = -1

# Row 144: restore the correct fixture data for this row
$ws.Range("B144").Value = 7532413
$ws.Range("F144").Value = 'Libertad Gran Mamore FC'
$ws.Range("G144").Value = 'Club Aurora'
$ws.Range("H144").Value = 0
$ws.Range("I144").Value = 1
$ws.Range("J144").Value = 'A'
$ws.Range("K144").Value = 2.25
$ws.Range("L144").Value = 3.3
$ws.Range("M144").Value = 2.8
$ws.Range("N144").Value = 2.375
$ws.Range("O144").Value = 3.4
$ws.Range("P144").Value = 2.875
$ws.Range("Q144").Value = -0.25
$ws.Range("R144").Value = 2.025
$ws.Range("S144").Value = 1.775
$ws.Range("T144").Value = 2.5
$ws.Range("U144").Value = 1.9
$ws.Range("V144").Value = 1.9
$ws.Range("W144").Value = -1
$ws.Range("X144").Value = -1
$ws.Range("Y144").Value = 1.875
$ws.Range("Z144").Value = -1
$ws.Range("AA144").Value = 0.7749999999999999
$ws.Range("AB144").Value = -1
$ws.Range("AC144").Value = 0.8999999999999999

# Row 145: restore the correct fixture data for this row
$ws.Range("B145").Value = 7532414
$ws.Range("F145").Value = 'Independiente Petrolero'
$ws.Range("G145").Value = 'Real Santa Cruz'
$ws.Range("H145").Value = 1
$ws.Range("I145").Value = 0
$ws.Range("J145").Value = 'H'
$ws.Range("K145").Value = 1.571
$ws.Range("L145").Value = 3.75
$ws.Range("M145").Value = 5
$ws.Range("N145").Value = 1.3
$ws.Range("O145").Value = 5
$ws.Range("P145").Value = 11
$ws.Range("Q145").Value = -1.75
$ws.Range("R145").Value = 2
$ws.Range("S145").Value = 1.8
$ws.Range("T145").Value = 3
$ws.Range("U145").Value = 1.85
$ws.Range("V145").Value = 1.95
$ws.Range("W145").Value = 0.3
$ws.Range("X145").Value = -1
$ws.Range("Y145").Value = -1
$ws.Range("Z145").Value = -1
$ws.Range("AA145").Value = 0.8
$ws.Range("AB145").Value = -1
$ws.Range("AC145").Value = 0.95

# Row 149: restore the correct fixture data for this row
$ws.Range("B149").Value = 7532420
$ws.Range("F149").Value = 'Club Aurora'
$ws.Range("G149").Value = 'Vaca Diez'
$ws.Range("H149").Value = 3
$ws.Range("I149").Value = 0
$ws.Range("J149").Value = 'H'
$ws.Range("K149").Value = 1.333
$ws.Range("L149").Value = 5
$ws.Range("M149").Value = 8
$ws.Range("N149").Value = 1.3
$ws.Range("O149").Value = 6.5
$ws.Range("P149").Value = 7
$ws.Range("Q149").Value = -1.5
$ws.Range("R149").Value = 1.8
$ws.Range("S149").Value = 2
$ws.Range("T149").Value = 3.25
$ws.Range("U149").Value = 1.95
$ws.Range("V149").Value = 1.85
$ws.Range("W149").Value = 0.3
$ws.Range("X149").Value = -1
$ws.Range("Y149").Value = -1
$ws.Range("Z149").Value = 0.8
$ws.Range("AA149").Value = -1
$ws.Range("AB149").Value = -0.5
$ws.Range("AC149").Value = 0.425

# Row 150: restore the correct fixture data for this row
$ws.Range("B150").Value = 7532421
$ws.Range("F150").Value = 'Guabira'
$ws.Range("G150").Value = 'Independiente Petrolero'
$ws.Range("H150").Value = 2
$ws.Range("I150").Value = 0
$ws.Range("J150").Value = 'H'
$ws.Range("K150").Value = 1.4
$ws.Range("L150").Value = 4.5
$ws.Range("M150").Value = 7.5
$ws.Range("N150").Value = 1.333
$ws.Range("O150").Value = 5.5
$ws.Range("P150").Value = 9.5
$ws.Range("Q150").Value = -1.5
$ws.Range("R150").Value = 1.85
$ws.Range("S150").Value = 1.95
$ws.Range("T150").Value = 3
$ws.Range("U150").Value = 1.825
$ws.Range("V150").Value = 1.975
$ws.Range("W150").Value = 0.333
$ws.Range("X150").Value = -1
$ws.Range("Y150").Value = -1
$ws.Range("Z150").Value = 0.8500000000000001
$ws.Range("AA150").Value = -1
$ws.Range("AB150").Value = -1
$ws.Range("AC150").Value = 0.9750000000000001

# Row 214: restore the correct fixture data for this row
$ws.Range("B214").Value = 8039392
$ws.Range("F214").Value = 'Oriente Petrolero'
$ws.Range("G214").Value = 'Jorge Wilstermann'
$ws.Range("H214").Value = 2
$ws.Range("I214").Value = 1
$ws.Range("J214").Value = 'H'
$ws.Range("K214").Value = 2
$ws.Range("L214").Value = 3.25
$ws.Range("M214").Value = 3.4
$ws.Range("N214").Value = 1.727
$ws.Range("O214").Value = 4
$ws.Range("P214").Value = 4.5
$ws.Range("Q214").Value = -0.75
$ws.Range("R214").Value = 1.9
$ws.Range("S214").Value = 1.9
$ws.Range("T214").Value = 2.75
$ws.Range("U214").Value = 1.9
$ws.Range("V214").Value = 1.9
$ws.Range("W214").Value = 0.7270000000000001
$ws.Range("X214").Value = -1
$ws.Range("Y214").Value = -1
$ws.Range("Z214").Value = 0.45
$ws.Range("AA214").Value = -0.5
$ws.Range("AB214").Value = 0.45
$ws.Range("AC214").Value = -0.5

# Row 215: restore the correct fixture data for this row
$ws.Range("B215").Value = 8038943
$ws.Range("F215").Value = 'San Jose de Oruro'
$ws.Range("G215").Value = 'Bolivar'
$ws.Range("H215").Value = 2
$ws.Range("I215").Value = 1
$ws.Range("J215").Value = 'H'
$ws.Range("K215").Value = 2.3
$ws.Range("L215").Value = 3.5
$ws.Range("M215").Value = 2.625
$ws.Range("N215").Value = 2.8
$ws.Range("O215").Value = 3.6
$ws.Range("P215").Value = 2.375
$ws.Range("Q215").Value = 0.25
$ws.Range("R215").Value = 1.8
$ws.Range("S215").Value = 2
$ws.Range("T215").Value = 3.25
$ws.Range("U215").Value = 1.975
$ws.Range("V215").Value = 1.825
$ws.Range("W215").Value = 1.8
$ws.Range("X215").Value = -1
$ws.Range("Y215").Value = -1
$ws.Range("Z215").Value = 0.8
$ws.Range("AA215").Value = -1
$ws.Range("AB215").Value = -0.5
$ws.Range("AC215").Value = 0.4125

# Remove the extraneous not-yet-played fixture row that trailed the sheet
# (id 216, match 8021846); dimension shrinks from AC218 to AC217.
$ws.Rows(218).Delete()
